# Update Orders.xlsx from orders.jsonl: append new order rows 3 and 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text (not get reinterpreted by
# Excel as a number/date), without leaving a residual cell style behind.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.NumberFormat = "General"
    $rng.ClearFormats()
}

# --- Row 3: ORD-2025002 ---
$ws.Range("A3").Value = "ORD-2025002"
$ws.Range("B3").Value = "فاطمة محمود"
Set-TextValue "C3" "0559876543"
$ws.Range("D3").Value = "أبو ظبي"
$ws.Range("E3").Value = "عطر فرنسي + شمعات عطرة"
$ws.Range("F3").Value = "350 AED"
Set-TextValue "G3" "2025-12-13"
$ws.Range("H3").Value = "2025-12-13T17:10:00Z"

# --- Row 4: ORD-2025003 ---
$ws.Range("A4").Value = "ORD-2025003"
$ws.Range("B4").Value = "علي محمد"
Set-TextValue "C4" "0507654321"
$ws.Range("D4").Value = "الشارقة"
$ws.Range("E4").Value = "ساعة ذهبية + بطاقة u03hf"
$ws.Range("F4").Value = "500 AED"
Set-TextValue "G4" "2025-12-13"
$ws.Range("H4").Value = "2025-12-13T17:20:00Z"
